$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(316).Insert()

$ws.Cells.Item(316, 1).Value = 5
$ws.Cells.Item(316, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(316, 3).Value = "Maule"
$ws.Cells.Item(316, 4).Value = 44900
$ws.Cells.Item(316, 5).Value = 7
$ws.Cells.Item(316, 6).Value = 100114014
$ws.Cells.Item(316, 7).Value = "Betarraga"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 5000
$ws.Cells.Item(316, 11).Value = 700
$ws.Cells.Item(316, 12).Value = 700
$ws.Cells.Item(316, 13).Value = 700
$ws.Cells.Item(316, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(316, 15).Value = "Región del Maule"
$ws.Cells.Item(316, 16).Value = 140
$ws.Cells.Item(316, 17).Value = 5
$ws.Cells.Item(316, 18).Value = "Hortaliza"
